{"js": "async (context) => {\n  // Map of old text -> new text, taken directly from the diff.\n  const replacements = [\n    [\"2024-05-20 Monday\", \"2024-05-21 Tuesday\"],\n    [\"51\u00f78=6, 3\", \"35\u00f73=11, 2\"],\n    [\"16\u00f73=5, 1\", \"75\u00f72=37, 1\"],\n    [\"36\u00f76=6, 0\", \"72\u00f74=18, 0\"],\n    [\"86\u00f72=43, 0\", \"51\u00f74=12, 3\"],\n    [\"52\u00f76=8, 4\", \"83\u00f72=41, 1\"],\n    [\"83\u00f78=10, 3\", \"90\u00f77=12, 6\"],\n    [\"55\u00f77=7, 6\", \"48\u00f74=12, 0\"],\n    [\"68\u00f76=11, 2\", \"89\u00f78=11, 1\"],\n    [\"66\u00f78=8, 2\", \"40\u00f77=5, 5\"],\n    [\"38\u00f78=4, 6\", \"31\u00f76=5, 1\"],\n    [\"91\u00f73=30, 1\", \"15\u00f75=3, 0\"],\n    [\"17\u00f75=3, 2\", \"90\u00f78=11, 2\"],\n    [\"95\u00f78=11, 7\", \"44\u00f77=6, 2\"],\n    [\"21\u00f74=5, 1\", \"48\u00f74=12, 0\"],\n    [\"71\u00f77=10, 1\", \"76\u00f74=19, 0\"],\n    [\"13\u00f73=4, 1\", \"24\u00f74=6, 0\"],\n    [\"25\u00f79=2, 7\", \"44\u00f73=14, 2\"],\n    [\"74\u00f72=37, 0\", \"88\u00f72=44, 0\"],\n    [\"98\u00f74=24, 2\", \"93\u00f73=31, 0\"],\n    [\"21\u00f75=4, 1\", \"45\u00f76=7, 3\"],\n    [\"49\u00f79=5, 4\", \"39\u00f75=7, 4\"],\n    [\"71\u00f78=8, 7\", \"72\u00f74=18, 0\"],\n    [\"99\u00f73=33, 0\", \"53\u00f77=7, 4\"],\n    [\"85\u00f75=17, 0\", \"55\u00f72=27, 1\"],\n    [\"70\u00f78=8, 6\", \"75\u00f73=25, 0\"],\n  ];\n\n  const body = context.document.body;\n\n  for (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (const range of results.items) {\n      range.insertText(newText, \"Replace\");\n    }\n    await context.sync();\n  }\n};\n", "ps1": "$d = $word.ActiveDocument\n\n# Old -> new text pairs, taken directly from the diff. All old values are\n# unique in the document, so a plain Find/Replace-all per pair is exact.\n$pairs = @(\n    @{ old = \"2024-05-20 Monday\"; new = \"2024-05-21 Tuesday\" },\n    @{ old = \"51\u00f78=6, 3\";  new = \"35\u00f73=11, 2\" },\n    @{ old = \"16\u00f73=5, 1\";  new = \"75\u00f72=37, 1\" },\n    @{ old = \"36\u00f76=6, 0\";  new = \"72\u00f74=18, 0\" },\n    @{ old = \"86\u00f72=43, 0\"; new = \"51\u00f74=12, 3\" },\n    @{ old = \"52\u00f76=8, 4\";  new = \"83\u00f72=41, 1\" },\n    @{ old = \"83\u00f78=10, 3\"; new = \"90\u00f77=12, 6\" },\n    @{ old = \"55\u00f77=7, 6\";  new = \"48\u00f74=12, 0\" },\n    @{ old = \"68\u00f76=11, 2\"; new = \"89\u00f78=11, 1\" },\n    @{ old = \"66\u00f78=8, 2\";  new = \"40\u00f77=5, 5\" },\n    @{ old = \"38\u00f78=4, 6\";  new = \"31\u00f76=5, 1\" },\n    @{ old = \"91\u00f73=30, 1\"; new = \"15\u00f75=3, 0\" },\n    @{ old = \"17\u00f75=3, 2\";  new = \"90\u00f78=11, 2\" },\n    @{ old = \"95\u00f78=11, 7\"; new = \"44\u00f77=6, 2\" },\n    @{ old = \"21\u00f74=5, 1\";  new = \"48\u00f74=12, 0\" },\n    @{ old = \"71\u00f77=10, 1\"; new = \"76\u00f74=19, 0\" },\n    @{ old = \"13\u00f73=4, 1\";  new = \"24\u00f74=6, 0\" },\n    @{ old = \"25\u00f79=2, 7\";  new = \"44\u00f73=14, 2\" },\n    @{ old = \"74\u00f72=37, 0\"; new = \"88\u00f72=44, 0\" },\n    @{ old = \"98\u00f74=24, 2\"; new = \"93\u00f73=31, 0\" },\n    @{ old = \"21\u00f75=4, 1\";  new = \"45\u00f76=7, 3\" },\n    @{ old = \"49\u00f79=5, 4\";  new = \"39\u00f75=7, 4\" },\n    @{ old = \"71\u00f78=8, 7\";  new = \"72\u00f74=18, 0\" },\n    @{ old = \"99\u00f73=33, 0\"; new = \"53\u00f77=7, 4\" },\n    @{ old = \"85\u00f75=17, 0\"; new = \"55\u00f72=27, 1\" },\n    @{ old = \"70\u00f78=8, 6\";  new = \"75\u00f73=25, 0\" }\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.old\n    $find.Replacement.Text = $pair.new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n"}
